# Generate Report for Handoff
# Update the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" values
# for the ca03050d-a2be-45e6-8869-2411d5b68e55.md file across all three sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Row 5 on every sheet corresponds to ca03050d-a2be-45e6-8869-2411d5b68e55.md
$wsOverview.Range("G5").Value = "2016-09-06 04:03:59"
$wsZhCn.Range("H5").Value     = "2016-09-06 04:03:46"
$wsDeDe.Range("H5").Value     = "2016-09-06 04:03:59"
